$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the report title (B1)
$ws.Range("B1").Value = "Atualização dos Indicadores de Confiança e de Clima Económico: Relatório da Ultima actividade"

# Update the job Name and description (C4, D4)
$ws.Range("C4").Value = "Aggregate Business Confidence"
$ws.Range("D4").Value = "Aggregate Business Confidence update"

# Update the run timestamp (G4)
$ws.Range("G4").Value = 44831.67255092743
